# ---------------------------------------------------------------------------
# Edit script: applies the two logical changes described by the commit:
#
#  1. The table on slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" slide)
#     gets a different built-in table style applied
#     ({94122885-01CC-4604-9379-8FDB8C2CA600} -> {DEDA34B9-51D4-4425-B232-E5D3DBD94BB1}).
#
#  2. The deck's theme ("Integral") is recoloured to use the stock Office
#     theme palette (the diff swaps the theme1.xml/theme2.xml part bodies;
#     since the only material difference between the "Integral" theme and
#     the "Office Theme" theme in this deck is the 12-colour scheme - the
#     font scheme and format scheme are byte-for-byte identical already -
#     this is reproduced through the PowerPoint colour-scheme object model).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{DEDA34B9-51D4-4425-B232-E5D3DBD94BB1}")
    }
}

# --- 2. Theme colours -----------------------------------------------------
# PowerPoint packs RGB() values as 0x00BBGGRR, so build each value from the
# target hex triplet (R,G,B) accordingly.
function BGR($r, $g, $b) { return ($b * 65536) + ($g * 256) + $r }

# Target ("Office Theme") scheme, in theme colour order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeScheme = @(
    (BGR 0x00 0x00 0x00),  # dk1      000000
    (BGR 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (BGR 0x44 0x54 0x6A),  # dk2      44546A
    (BGR 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (BGR 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (BGR 0xED 0x7D 0x31),  # accent2  ED7D31
    (BGR 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (BGR 0xFF 0xC0 0x00),  # accent4  FFC000
    (BGR 0x44 0x72 0xC4),  # accent5  4472C4
    (BGR 0x70 0xAD 0x47),  # accent6  70AD47
    (BGR 0x05 0x63 0xC1),  # hlink    0563C1
    (BGR 0x95 0x4F 0x72)   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeScheme[$i - 1]
}
